$d = $word.ActiveDocument

$d.Content.Find.Execute("Fall 2022 INFO-233", $true, $false, $false, $false, $false,
                         $true, 1, $false, "INFO-233", 2)
